$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the error-description cells in row 9
$ws.Range("B9").Value = "Koszyk i teczki akt osobowych"
$ws.Range("C9").Value = "Błędne komunikaty i brak nazwiska i imienia"

# Row 9 grew taller to fit the new text
$ws.Rows.Item(9).RowHeight = 42.75

# Move the active selection to D9
$ws.Range("D9").Select()
